# Actualizacion automatica hashcode - update MD5 hash values in column B
# (column A = item code, column B = hashcode) for the rows whose hash changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

if ($ws.Range("A9").Value2 -ne "05-050305TC") { throw "Unexpected row 9 content" }
$ws.Range("B9").Value = "25caa2dbb6863594665f9347254df2eb"

if ($ws.Range("A11").Value2 -ne "05-050301A") { throw "Unexpected row 11 content" }
$ws.Range("B11").Value = "24ba6a91b9f310c879a707c8548bc5ca"

if ($ws.Range("A17").Value2 -ne "05-050305TP") { throw "Unexpected row 17 content" }
$ws.Range("B17").Value = "f84965b7eda5867e2ee4214735148b8a"

if ($ws.Range("A34").Value2 -ne "05-050316TP") { throw "Unexpected row 34 content" }
$ws.Range("B34").Value = "f26b8661c6953e87a1e135d1ce4ba4f9"

if ($ws.Range("A126").Value2 -ne "05-050309A") { throw "Unexpected row 126 content" }
$ws.Range("B126").Value = "5440b946acc7c09e85d1ea2b12fda6b3"

if ($ws.Range("A133").Value2 -ne "05-050312TP") { throw "Unexpected row 133 content" }
$ws.Range("B133").Value = "4409143d57b5150097d5d36c17aa15f5"

if ($ws.Range("A136").Value2 -ne "05-050312TC") { throw "Unexpected row 136 content" }
$ws.Range("B136").Value = "84b9fadbbe728ba09cb294c84cf1cc6c"

if ($ws.Range("A160").Value2 -ne "05-050203TP") { throw "Unexpected row 160 content" }
$ws.Range("B160").Value = "541092bd7be88459577967d0c5c185d6"

if ($ws.Range("A163").Value2 -ne "05-050308A") { throw "Unexpected row 163 content" }
$ws.Range("B163").Value = "d02109d78d059c69e67930e83c9ddf01"

if ($ws.Range("A170").Value2 -ne "05-050203TC") { throw "Unexpected row 170 content" }
$ws.Range("B170").Value = "ecda67c35968722c0c06df0d3bfd65ab"

if ($ws.Range("A176").Value2 -ne "05-050303TP") { throw "Unexpected row 176 content" }
$ws.Range("B176").Value = "b1f08c9b9593a4bbd64f3c68b0086dba"

if ($ws.Range("A181").Value2 -ne "05-050303TC") { throw "Unexpected row 181 content" }
$ws.Range("B181").Value = "803a55a9f4255f6dc2a4d211ac6630fd"

if ($ws.Range("A184").Value2 -ne "05-050305A") { throw "Unexpected row 184 content" }
$ws.Range("B184").Value = "a7a4fdebab5449d7b192a50e15e7d6c9"

if ($ws.Range("A192").Value2 -ne "05-050314TP") { throw "Unexpected row 192 content" }
$ws.Range("B192").Value = "3969bb9ea333d1f3d19157823fe04a57"

if ($ws.Range("A199").Value2 -ne "05-050314TC") { throw "Unexpected row 199 content" }
$ws.Range("B199").Value = "e2cd9281650b2561cce6e981c5071842"

if ($ws.Range("A201").Value2 -ne "05-050306A") { throw "Unexpected row 201 content" }
$ws.Range("B201").Value = "23ab0c1aebe30df9f876b2802b58f35e"

if ($ws.Range("A214").Value2 -ne "05-050303A") { throw "Unexpected row 214 content" }
$ws.Range("B214").Value = "d0871aa81a228cdf44e3daa7125e033b"

if ($ws.Range("A229").Value2 -ne "05-050205TP") { throw "Unexpected row 229 content" }
$ws.Range("B229").Value = "9d28ee6f80d1fb9989e4de321fbe309f"

if ($ws.Range("A230").Value2 -ne "05-050304A") { throw "Unexpected row 230 content" }
$ws.Range("B230").Value = "3639bf317e3adcaf83ba1a7259559e15"

if ($ws.Range("A234").Value2 -ne "05-050205TC") { throw "Unexpected row 234 content" }
$ws.Range("B234").Value = "b93283fe02b1652083f00ca9a4863e8b"

if ($ws.Range("A287").Value2 -ne "05-050201TC") { throw "Unexpected row 287 content" }
$ws.Range("B287").Value = "4239d5d0c4fcf7aca447448ae4041393"

if ($ws.Range("A299").Value2 -ne "05-050310TC") { throw "Unexpected row 299 content" }
$ws.Range("B299").Value = "dc79d5cdd9556ba6b68d5fd801d4df53"

if ($ws.Range("A308").Value2 -ne "05-050310TP") { throw "Unexpected row 308 content" }
$ws.Range("B308").Value = "202d47ae58a1e147e3e945c5cfcd9089"

if ($ws.Range("A345").Value2 -ne "05-050201TP") { throw "Unexpected row 345 content" }
$ws.Range("B345").Value = "a711563c8e65422475895cb04bc9c44e"

if ($ws.Range("A470").Value2 -ne "05-050204A") { throw "Unexpected row 470 content" }
$ws.Range("B470").Value = "bb6885cbc3f82f8ebb6168dbd13ab969"

if ($ws.Range("A489").Value2 -ne "05-050205A") { throw "Unexpected row 489 content" }
$ws.Range("B489").Value = "c7c23d5fcf4008c27d5a7a4435b78a3e"

if ($ws.Range("A491").Value2 -ne "05-050314A") { throw "Unexpected row 491 content" }
$ws.Range("B491").Value = "e14fe01c910387baaad5cba5ac23c98e"

if ($ws.Range("A505").Value2 -ne "05-050208TC") { throw "Unexpected row 505 content" }
$ws.Range("B505").Value = "666443888eb954a6ec192d1d2dc9846b"

if ($ws.Range("A514").Value2 -ne "05-050311A") { throw "Unexpected row 514 content" }
$ws.Range("B514").Value = "f1dfd8a0a27197abe3d01a6eb87d2e59"

if ($ws.Range("A520").Value2 -ne "05-050306TP") { throw "Unexpected row 520 content" }
$ws.Range("B520").Value = "4675c67bf2a8dc16775ec05abb7d5af3"

if ($ws.Range("A528").Value2 -ne "05-050317TC") { throw "Unexpected row 528 content" }
$ws.Range("B528").Value = "0c6ae3d60f548d9aaf15ba28c7ac77f9"

if ($ws.Range("A529").Value2 -ne "05-050312A") { throw "Unexpected row 529 content" }
$ws.Range("B529").Value = "7febf5349f4310f03801db71426221db"

if ($ws.Range("A531").Value2 -ne "05-050203A") { throw "Unexpected row 531 content" }
$ws.Range("B531").Value = "4f127ed49f06fad3e2353920c2c3ec3b"

if ($ws.Range("A539").Value2 -ne "05-050317TP") { throw "Unexpected row 539 content" }
$ws.Range("B539").Value = "927fdd666ff5c2131184c7611ca11117"

if ($ws.Range("A563").Value2 -ne "05-050201A") { throw "Unexpected row 563 content" }
$ws.Range("B563").Value = "b012fa5019c809c8b5b7d08521b80490"

if ($ws.Range("A566").Value2 -ne "05-050310A") { throw "Unexpected row 566 content" }
$ws.Range("B566").Value = "1e3f6acbfec19073c1e511d513db495d"

if ($ws.Range("A579").Value2 -ne "05-050308TC") { throw "Unexpected row 579 content" }
$ws.Range("B579").Value = "27fd39202f9a753ff30191dc7b415b9a"

if ($ws.Range("A588").Value2 -ne "05-050308TP") { throw "Unexpected row 588 content" }
$ws.Range("B588").Value = "921f91ea685cd04e378d975ff64e1987"

if ($ws.Range("A632").Value2 -ne "05-050204TP") { throw "Unexpected row 632 content" }
$ws.Range("B632").Value = "db677500f068667c02cca9c73c7ad109"

if ($ws.Range("A643").Value2 -ne "05-050204TC") { throw "Unexpected row 643 content" }
$ws.Range("B643").Value = "9c492c3bf149444381b4948a101b2855"

if ($ws.Range("A651").Value2 -ne "05-050302TC") { throw "Unexpected row 651 content" }
$ws.Range("B651").Value = "cff54a9e4d0702d29363119765df9c28"

if ($ws.Range("A671").Value2 -ne "05-050313TC") { throw "Unexpected row 671 content" }
$ws.Range("B671").Value = "d641df052017d6be477a892d14ba531a"

if ($ws.Range("A682").Value2 -ne "05-050317A") { throw "Unexpected row 682 content" }
$ws.Range("B682").Value = "50d6b3928ae51952c9f11b33b97961e0"

if ($ws.Range("A696").Value2 -ne "05-050206TP") { throw "Unexpected row 696 content" }
$ws.Range("B696").Value = "ce3302a18b89e77538c792d1acdf12b7"

if ($ws.Range("A701").Value2 -ne "05-050206TC") { throw "Unexpected row 701 content" }
$ws.Range("B701").Value = "24f934f507a6461a5b7d6a9df32b0683"

if ($ws.Range("A716").Value2 -ne "05-050304TC") { throw "Unexpected row 716 content" }
$ws.Range("B716").Value = "c6181dd5979b5d646a95ca80e3c0611b"

if ($ws.Range("A719").Value2 -ne "05-050206A") { throw "Unexpected row 719 content" }
$ws.Range("B719").Value = "104badc673a565577d2d251e733eeddb"

if ($ws.Range("A720").Value2 -ne "05-050315A") { throw "Unexpected row 720 content" }
$ws.Range("B720").Value = "194d96116d0b83bc7346b5f030d7ef73"

if ($ws.Range("A731").Value2 -ne "05-050304TP") { throw "Unexpected row 731 content" }
$ws.Range("B731").Value = "f518377aebbc9299f15cc759ab235084"

if ($ws.Range("A745").Value2 -ne "05-050316A") { throw "Unexpected row 745 content" }
$ws.Range("B745").Value = "0867eed9183bdebf6cc2ae2672c200c2"

if ($ws.Range("A758").Value2 -ne "05-050315TP") { throw "Unexpected row 758 content" }
$ws.Range("B758").Value = "a5326aa5e29f014ac571870f665bb88d"

if ($ws.Range("A845").Value2 -ne "05-050202TC") { throw "Unexpected row 845 content" }
$ws.Range("B845").Value = "81ac0b59df2a829838602f6dd1f016f8"

if ($ws.Range("A848").Value2 -ne "05-050311TC") { throw "Unexpected row 848 content" }
$ws.Range("B848").Value = "9310778aa61db567cecbf2a5058225ce"

if ($ws.Range("A853").Value2 -ne "05-050311TP") { throw "Unexpected row 853 content" }
$ws.Range("B853").Value = "a4ae71c034983d667a5049453166787f"

if ($ws.Range("A880").Value2 -ne "05-050309TP") { throw "Unexpected row 880 content" }
$ws.Range("B880").Value = "4d5d99a576159a067689719f6e6f6e99"
